$d = $word.ActiveDocument
$wNs = 'http://schemas.openxmlformats.org/wordprocessingml/2006/main'

# Capture stable paragraph references up-front (by index, 1-based) before any
# structural edits happen.  Object references keep working even after other
# paragraphs' contents are replaced via InsertXML.
$pPcAddr      = $d.Paragraphs(2)   # "1)   PC is initialized to start at 96.  Check cache for addr."
$pUseMasks    = $d.Paragraphs(4)   # "b)   Use masks to grab:  setIndex, tag ..."
$pGoToSet     = $d.Paragraphs(5)   # "c)  Go to correct set within cache (via setIndex). ..."
$pIfMatch     = $d.Paragraphs(6)   # "i)  If match (V is 1 AND tag match)"
$pLoadInstr   = $d.Paragraphs(7)   # "1)  Load instruction to preIssue buffer."
$pIfNotMatch  = $d.Paragraphs(10)  # "ii) If not match: (V is 0, OR ... no tag match)"
$pIncPcLast   = $d.Paragraphs(15)  # "5)  Inc PC." (last numbered sub-step before trailing blank paragraph)

$tabsPPr = '<w:pPr><w:tabs><w:tab w:val="left" w:pos="360"/><w:tab w:val="left" w:pos="720"/><w:tab w:val="left" w:pos="1080"/><w:tab w:val="left" w:pos="1440"/></w:tabs></w:pPr>'

# --- Change 1: "Check cache for addr." -> wrap "addr" in proofErr spell tags ---
$xml1 = '<w:p xmlns:w="' + $wNs + '">' + $tabsPPr + `
  '<w:r><w:t xml:space="preserve">1)   PC is initialized to start at 96.  Check cache for </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:t>addr</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t>.</w:t></w:r>' + `
  '</w:p>'
$pPcAddr.Range.InsertXML($xml1)

# --- Change 2: "Use masks to grab:  setIndex, tag" -> wrap "setIndex" in proofErr spell tags ---
$xml2 = '<w:p xmlns:w="' + $wNs + '">' + $tabsPPr + `
  '<w:r><w:tab/><w:t xml:space="preserve">b)   Use masks to grab:  </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:t>setIndex</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t>, tag</w:t></w:r>' + `
  '<w:r w:rsidR="00466F7F"><w:t xml:space="preserve"> (no need for byte offset for address)</w:t></w:r>' + `
  '<w:r><w:t>.</w:t></w:r>' + `
  '</w:p>'
$pUseMasks.Range.InsertXML($xml2)

# --- Change 3: "Go to correct set within cache (via setIndex)." -> wrap "setIndex" in proofErr spell tags ---
$xml3 = '<w:p xmlns:w="' + $wNs + '">' + $tabsPPr + `
  '<w:r><w:tab/><w:t xml:space="preserve">c)  Go to correct set within cache (via </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:t>setIndex</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t>).  Test each block for tag match.</w:t></w:r>' + `
  '</w:p>'
$pGoToSet.Range.InsertXML($xml3)

# --- Change 4: "i)  If match" -> wrap "i" in proofErr spell tags ---
$xml4 = '<w:p xmlns:w="' + $wNs + '">' + $tabsPPr + `
  '<w:r><w:tab/></w:r>' + `
  '<w:r><w:tab/></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:t>i</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t>)  If match</w:t></w:r>' + `
  '<w:r w:rsidR="000F7FC4"><w:t xml:space="preserve"> (V is 1 AND tag match)</w:t></w:r>' + `
  '</w:p>'
$pIfMatch.Range.InsertXML($xml4)

# --- Change 5: "Load instruction to preIssue buffer." -> wrap "preIssue" in proofErr spell tags ---
$xml5 = '<w:p xmlns:w="' + $wNs + '">' + $tabsPPr + `
  '<w:r><w:tab/></w:r>' + `
  '<w:r><w:tab/></w:r>' + `
  '<w:r><w:tab/><w:t xml:space="preserve">1)  Load instruction to </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:t>preIssue</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> buffer.</w:t></w:r>' + `
  '</w:p>'
$pLoadInstr.Range.InsertXML($xml5)

# --- Change 6: "ii) If not match: (V is 0, OR no tag match)" -> merge into a single run,
#     drop the _GoBack bookmark that used to sit mid-sentence here ---
$xml6 = '<w:p xmlns:w="' + $wNs + '">' + $tabsPPr + `
  '<w:r><w:tab/></w:r>' + `
  '<w:r><w:tab/><w:t xml:space="preserve">ii) If not match: </w:t></w:r>' + `
  '<w:r w:rsidR="000F7FC4"><w:t>(V is 0, OR no tag match)</w:t></w:r>' + `
  '</w:p>'
$pIfNotMatch.Range.InsertXML($xml6)

# --- Change 7: append a brand-new paragraph after "5)  Inc PC." holding the
#     final remark and the relocated _GoBack bookmark ---
$newRange = $pIncPcLast.Range.InsertParagraphAfter()
$newPara = $pIncPcLast.Next()
$xml7 = '<w:p xmlns:w="' + $wNs + '">' + $tabsPPr + `
  '<w:r><w:t>Only data memory is in cache?????</w:t></w:r>' + `
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
  '<w:bookmarkEnd w:id="0"/>' + `
  '</w:p>'
$newPara.Range.InsertXML($xml7)
